$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 04:12"

# Update the "countries / provincias" block (rows 99-108): the country
# order shifted (Ghana moved down, after "Estado de Palestina") and the
# case numbers for each row were refreshed.
$ws.Range("A99").Value = "Azerbaiyan"
$ws.Range("B99").Value = 93
$ws.Range("C99").Value = 6
$ws.Range("D99").Value = 10
$ws.Range("E99").Value = 81
$ws.Range("F99").Value = 6
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 2

$ws.Range("A100").Value = "Bielorrusia"
$ws.Range("B100").Value = 86
$ws.Range("C100").Value = 5
$ws.Range("D100").Value = 29
$ws.Range("E100").Value = 57
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0

$ws.Range("A101").Value = "Afganistan"
$ws.Range("B101").Value = 84
$ws.Range("C101").Value = 10
$ws.Range("D101").Value = 2
$ws.Range("E101").Value = 80
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 2

$ws.Range("A102").Value = "Kazajistan"
$ws.Range("B102").Value = 81
$ws.Range("C102").Value = 9
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 81
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0

$ws.Range("A103").Value = "Costa de Marfil"
$ws.Range("B103").Value = 80
$ws.Range("C103").Value = 7
$ws.Range("D103").Value = 3
$ws.Range("E103").Value = 77
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

$ws.Range("A104").Value = "Camerun"
$ws.Range("B104").Value = 75
$ws.Range("C104").Value = 9
$ws.Range("D104").Value = 2
$ws.Range("E104").Value = 72
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 1

$ws.Range("A105").Value = "Georgia"
$ws.Range("B105").Value = 75
$ws.Range("C105").Value = 5
$ws.Range("D105").Value = 10
$ws.Range("E105").Value = 65
$ws.Range("F105").Value = 1
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

$ws.Range("A106").Value = "Guadalupe"
$ws.Range("B106").Value = 73
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 72
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 1

$ws.Range("A107").Value = "Estado de Palestina"
$ws.Range("B107").Value = 71
$ws.Range("C107").Value = 11
$ws.Range("D107").Value = 16
$ws.Range("E107").Value = 54
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 1

$ws.Range("A108").Value = "Ghana"
$ws.Range("B108").Value = 68
$ws.Range("C108").Value = 15
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 64
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 4
